$d = $word.ActiveDocument

$replacements = @(
    @{Old = "355×7=2485"; New = "361×3=1083"}
    @{Old = "417×6=2502"; New = "869×5=4345"}
    @{Old = "888×2=1776"; New = "235×9=2115"}
    @{Old = "445×3=1335"; New = "316×2=632"}
    @{Old = "899×2=1798"; New = "335×4=1340"}
    @{Old = "862×4=3448"; New = "560×3=1680"}
    @{Old = "540×7=3780"; New = "226×2=452"}
    @{Old = "995×5=4975"; New = "982×4=3928"}
    @{Old = "902×8=7216"; New = "747×5=3735"}
    @{Old = "761×9=6849"; New = "601×6=3606"}
    @{Old = "560×5=2800"; New = "145×8=1160"}
    @{Old = "347×4=1388"; New = "430×2=860"}
    @{Old = "939×5=4695"; New = "825×4=3300"}
    @{Old = "123×6=738"; New = "537×9=4833"}
    @{Old = "322×5=1610"; New = "795×9=7155"}
    @{Old = "462×9=4158"; New = "179×3=537"}
    @{Old = "245×8=1960"; New = "336×4=1344"}
    @{Old = "652×4=2608"; New = "769×4=3076"}
    @{Old = "839×7=5873"; New = "382×7=2674"}
    @{Old = "667×6=4002"; New = "954×3=2862"}
    @{Old = "103×6=618"; New = "971×2=1942"}
    @{Old = "709×4=2836"; New = "862×5=4310"}
    @{Old = "274×9=2466"; New = "191×2=382"}
    @{Old = "388×2=776"; New = "759×5=3795"}
    @{Old = "933×5=4665"; New = "376×8=3008"}
)

foreach ($r in $replacements) {
    $found = $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
    if (-not $found) {
        Write-Host "WARNING: replacement not found for $($r.Old)"
    }
}

$d.Save()
